$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 2000
$ws.Cells.Item(40, 10).Value = 0
$ws.Cells.Item(40, 12).Value = 0
$ws.Cells.Item(113, 8).Value = 82056.16
$ws.Cells.Item(113, 9).Value = 252433
$ws.Cells.Item(113, 10).Value = 6333.1113
$ws.Cells.Item(113, 11).Value = 252433
$ws.Cells.Item(113, 12).Value = 6333.1113
$ws.Cells.Item(113, 13).Value = -249179
$ws.Cells.Item(113, 14).Value = -12841.1113
$ws.Cells.Item(132, 8).Value = 2270.151
$ws.Cells.Item(132, 9).Value = 1261.0212
$ws.Cells.Item(132, 10).Value = 10175
$ws.Cells.Item(132, 11).Value = 3783.063599999999
$ws.Cells.Item(132, 12).Value = 30525
$ws.Cells.Item(132, 13).Value = -1253.063599999999
$ws.Cells.Item(132, 14).Value = -35585
$ws.Cells.Item(137, 8).Value = 1175.1428
$ws.Cells.Item(137, 9).Value = 1090.2727
$ws.Cells.Item(137, 11).Value = 3270.8181
$ws.Cells.Item(137, 13).Value = -720.8181
$ws.Cells.Item(141, 8).Value = 3050.5
$ws.Cells.Item(141, 9).Value = 3698.8
$ws.Cells.Item(141, 11).Value = 11096.4
$ws.Cells.Item(141, 13).Value = -5916.400000000001
$ws.Cells.Item(40, 14).ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 3566.0364
$ws.Cells.Item(32, 9).Value = 2986.7
$ws.Cells.Item(32, 10).Value = 9359.4
$ws.Cells.Item(32, 11).Value = 2986.7
$ws.Cells.Item(32, 12).Value = 9359.4
$ws.Cells.Item(32, 13).Value = -2699.7
$ws.Cells.Item(32, 14).Value = -9933.4
$ws.Cells.Item(61, 8).Value = 1559.5238
$ws.Cells.Item(61, 9).Value = 1153.4117
$ws.Cells.Item(61, 10).Value = 3285.5
$ws.Cells.Item(61, 11).Value = 1153.4117
$ws.Cells.Item(61, 12).Value = 3285.5
$ws.Cells.Item(61, 13).Value = -941.4117000000001
$ws.Cells.Item(61, 14).Value = -3709.5
$ws.Cells.Item(97, 8).Value = 3547.3333
$ws.Cells.Item(97, 9).Value = 3547.3333
$ws.Cells.Item(97, 11).Value = 3547.3333
$ws.Cells.Item(97, 13).Value = -3051.3333
$ws.Cells.Item(122, 8).Value = 2883.2307
$ws.Cells.Item(122, 9).Value = 2883.2307
$ws.Cells.Item(122, 11).Value = 8649.6921
$ws.Cells.Item(122, 13).Value = -6199.6921
$ws.Cells.Item(136, 8).Value = 1559.5238
$ws.Cells.Item(136, 9).Value = 1153.4117
$ws.Cells.Item(136, 10).Value = 3285.5
$ws.Cells.Item(136, 11).Value = 3460.2351
$ws.Cells.Item(136, 12).Value = 9856.5
$ws.Cells.Item(136, 13).Value = -910.2351000000003
$ws.Cells.Item(136, 14).Value = -14956.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 4661.7417
$ws.Cells.Item(20, 9).Value = 5178.7827
$ws.Cells.Item(20, 11).Value = 5178.7827
$ws.Cells.Item(20, 13).Value = -4931.7827
$ws.Cells.Item(99, 8).Value = 2750
$ws.Cells.Item(99, 9).Value = 2000
$ws.Cells.Item(99, 11).Value = 2000
$ws.Cells.Item(99, 13).Value = -502
$ws.Cells.Item(107, 8).Value = 92817.91
$ws.Cells.Item(107, 9).Value = 111888.664
$ws.Cells.Item(107, 10).Value = 6999.5
$ws.Cells.Item(107, 11).Value = 111888.664
$ws.Cells.Item(107, 12).Value = 6999.5
$ws.Cells.Item(107, 13).Value = -109968.664
$ws.Cells.Item(107, 14).Value = -10839.5
$ws.Cells.Item(134, 8).Value = 1358.8422
$ws.Cells.Item(134, 9).Value = 1362.8889
$ws.Cells.Item(134, 11).Value = 4088.6667
$ws.Cells.Item(134, 13).Value = -1553.6667

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1813.2222
$ws.Cells.Item(16, 10).Value = 2009.4
$ws.Cells.Item(16, 12).Value = 2009.4
$ws.Cells.Item(16, 14).Value = -2583.4
$ws.Cells.Item(31, 8).Value = 19689.684
$ws.Cells.Item(31, 9).Value = 2673.6667
$ws.Cells.Item(31, 10).Value = 27543.23
$ws.Cells.Item(31, 11).Value = 2673.6667
$ws.Cells.Item(31, 12).Value = 27543.23
$ws.Cells.Item(31, 13).Value = -2378.6667
$ws.Cells.Item(31, 14).Value = -28133.23
$ws.Cells.Item(34, 8).Value = 19689.684
$ws.Cells.Item(34, 9).Value = 2673.6667
$ws.Cells.Item(34, 10).Value = 27543.23
$ws.Cells.Item(34, 11).Value = 2673.6667
$ws.Cells.Item(34, 12).Value = 27543.23
$ws.Cells.Item(34, 13).Value = -2471.6667
$ws.Cells.Item(34, 14).Value = -27947.23
$ws.Cells.Item(38, 8).Value = 10347.556
$ws.Cells.Item(38, 9).Value = 3825.6
$ws.Cells.Item(38, 10).Value = 18500
$ws.Cells.Item(38, 11).Value = 3825.6
$ws.Cells.Item(38, 12).Value = 18500
$ws.Cells.Item(38, 13).Value = -3448.6
$ws.Cells.Item(38, 14).Value = -19254
$ws.Cells.Item(46, 8).Value = 10347.556
$ws.Cells.Item(46, 9).Value = 3825.6
$ws.Cells.Item(46, 10).Value = 18500
$ws.Cells.Item(46, 11).Value = 3825.6
$ws.Cells.Item(46, 12).Value = 18500
$ws.Cells.Item(46, 13).Value = -3614.6
$ws.Cells.Item(46, 14).Value = -18922
$ws.Cells.Item(58, 8).Value = 1866
$ws.Cells.Item(58, 9).Value = 1299.5
$ws.Cells.Item(58, 11).Value = 1299.5
$ws.Cells.Item(58, 13).Value = -1096.5
$ws.Cells.Item(113, 8).Value = 1813.2222
$ws.Cells.Item(113, 10).Value = 2009.4
$ws.Cells.Item(113, 12).Value = 2009.4
$ws.Cells.Item(113, 14).Value = -6349.4
$ws.Cells.Item(122, 8).Value = 200922.6
$ws.Cells.Item(122, 9).Value = 125875.625
$ws.Cells.Item(122, 11).Value = 377626.875
$ws.Cells.Item(122, 13).Value = -375176.875
$ws.Cells.Item(132, 8).Value = 2312.4348
$ws.Cells.Item(132, 9).Value = 2189.8572
$ws.Cells.Item(132, 10).Value = 3599.5
$ws.Cells.Item(132, 11).Value = 6569.571599999999
$ws.Cells.Item(132, 12).Value = 10798.5
$ws.Cells.Item(132, 13).Value = -4039.571599999999
$ws.Cells.Item(132, 14).Value = -15858.5
$ws.Cells.Item(134, 8).Value = 2282.6538
$ws.Cells.Item(134, 9).Value = 1276.0869
$ws.Cells.Item(134, 10).Value = 9999.666999999999
$ws.Cells.Item(134, 11).Value = 3828.2607
$ws.Cells.Item(134, 12).Value = 29999.001
$ws.Cells.Item(134, 13).Value = -1293.2607
$ws.Cells.Item(134, 14).Value = -35069.001
$ws.Cells.Item(136, 8).Value = 1866
$ws.Cells.Item(136, 9).Value = 1299.5
$ws.Cells.Item(136, 11).Value = 3898.5
$ws.Cells.Item(136, 13).Value = -1348.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 2076.4167
$ws.Cells.Item(5, 9).Value = 1058.7778
$ws.Cells.Item(5, 10).Value = 5129.3335
$ws.Cells.Item(5, 11).Value = 3176.3334
$ws.Cells.Item(5, 12).Value = 15388.0005
$ws.Cells.Item(5, 13).Value = -3064.3334
$ws.Cells.Item(5, 14).Value = -15612.0005
$ws.Cells.Item(25, 8).Value = 416.66666
$ws.Cells.Item(25, 9).Value = 250
$ws.Cells.Item(25, 10).Value = 500
$ws.Cells.Item(25, 11).Value = 750
$ws.Cells.Item(25, 12).Value = 1500
$ws.Cells.Item(25, 13).Value = -581
$ws.Cells.Item(25, 14).Value = -1838
$ws.Cells.Item(30, 8).Value = 416.66666
$ws.Cells.Item(30, 9).Value = 250
$ws.Cells.Item(30, 10).Value = 500
$ws.Cells.Item(30, 11).Value = 750
$ws.Cells.Item(30, 12).Value = 1500
$ws.Cells.Item(30, 13).Value = -648
$ws.Cells.Item(30, 14).Value = -1704
$ws.Cells.Item(63, 8).Value = 0
$ws.Cells.Item(63, 9).Value = 0
$ws.Cells.Item(63, 10).Value = 0
$ws.Cells.Item(63, 11).Value = 0
$ws.Cells.Item(63, 12).Value = 0
$ws.Cells.Item(66, 8).Value = 0
$ws.Cells.Item(66, 9).Value = 0
$ws.Cells.Item(66, 10).Value = 0
$ws.Cells.Item(66, 11).Value = 0
$ws.Cells.Item(66, 12).Value = 0
$ws.Cells.Item(107, 8).Value = 741.32355
$ws.Cells.Item(107, 9).Value = 1000
$ws.Cells.Item(107, 10).Value = 725.15625
$ws.Cells.Item(107, 11).Value = 3000
$ws.Cells.Item(107, 12).Value = 2175.46875
$ws.Cells.Item(107, 13).Value = -1080
$ws.Cells.Item(107, 14).Value = -6015.46875
$ws.Cells.Item(122, 8).Value = 1265.8
$ws.Cells.Item(122, 9).Value = 599
$ws.Cells.Item(122, 10).Value = 1432.5
$ws.Cells.Item(122, 11).Value = 5391
$ws.Cells.Item(122, 12).Value = 12892.5
$ws.Cells.Item(122, 13).Value = -2941
$ws.Cells.Item(122, 14).Value = -17792.5
$ws.Cells.Item(134, 8).Value = 5232.75
$ws.Cells.Item(134, 9).Value = 3166
$ws.Cells.Item(134, 10).Value = 8677.333000000001
$ws.Cells.Item(134, 11).Value = 9498
$ws.Cells.Item(134, 12).Value = 26031.999
$ws.Cells.Item(134, 13).Value = -4428
$ws.Cells.Item(134, 14).Value = -36171.999
$ws.Cells.Item(135, 8).Value = 2076.4167
$ws.Cells.Item(135, 9).Value = 1058.7778
$ws.Cells.Item(135, 10).Value = 5129.3335
$ws.Cells.Item(135, 11).Value = 9529.0002
$ws.Cells.Item(135, 12).Value = 46164.0015
$ws.Cells.Item(135, 13).Value = -6994.0002
$ws.Cells.Item(135, 14).Value = -51234.0015
$ws.Cells.Item(137, 8).Value = 5002388
$ws.Cells.Item(137, 9).Value = 10001549
$ws.Cells.Item(137, 10).Value = 3226.9
$ws.Cells.Item(137, 11).Value = 30004647
$ws.Cells.Item(137, 12).Value = 9680.700000000001
$ws.Cells.Item(137, 13).Value = -29999547
$ws.Cells.Item(137, 14).Value = -19880.7
$ws.Cells.Item(63, 13).ClearContents()
$ws.Cells.Item(63, 14).ClearContents()
$ws.Cells.Item(66, 13).ClearContents()
$ws.Cells.Item(66, 14).ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 22769.732
$ws.Cells.Item(97, 9).Value = 31710.238
$ws.Cells.Item(97, 11).Value = 31710.238
$ws.Cells.Item(97, 13).Value = -31214.238
$ws.Cells.Item(122, 8).Value = 1433.8235
$ws.Cells.Item(122, 9).Value = 1438.3846
$ws.Cells.Item(122, 10).Value = 1419
$ws.Cells.Item(122, 11).Value = 4315.1538
$ws.Cells.Item(122, 12).Value = 4257
$ws.Cells.Item(122, 13).Value = -1865.1538
$ws.Cells.Item(122, 14).Value = -9157
$ws.Cells.Item(132, 8).Value = 2207.6
$ws.Cells.Item(132, 9).Value = 2270.5833
$ws.Cells.Item(132, 10).Value = 696
$ws.Cells.Item(132, 11).Value = 6811.749899999999
$ws.Cells.Item(132, 12).Value = 2088
$ws.Cells.Item(132, 13).Value = -4281.749899999999
$ws.Cells.Item(132, 14).Value = -7148

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 48571.285
$ws.Cells.Item(7, 9).Value = 54999.832
$ws.Cells.Item(7, 11).Value = 54999.832
$ws.Cells.Item(7, 13).Value = -54887.832
$ws.Cells.Item(38, 8).Value = 10666.667
$ws.Cells.Item(38, 9).Value = 8000
$ws.Cells.Item(38, 11).Value = 8000
$ws.Cells.Item(38, 13).Value = -7590
$ws.Cells.Item(46, 8).Value = 26984.578
$ws.Cells.Item(46, 9).Value = 51466.555
$ws.Cells.Item(46, 10).Value = 4950.8
$ws.Cells.Item(46, 11).Value = 51466.555
$ws.Cells.Item(46, 12).Value = 4950.8
$ws.Cells.Item(46, 13).Value = -51278.555
$ws.Cells.Item(46, 14).Value = -5326.8
$ws.Cells.Item(93, 8).Value = 20292.334
$ws.Cells.Item(93, 9).Value = 1955.3334
$ws.Cells.Item(93, 11).Value = 1955.3334
$ws.Cells.Item(93, 13).Value = -707.3334
$ws.Cells.Item(122, 8).Value = 289214.44
$ws.Cells.Item(122, 9).Value = 336666.84
$ws.Cells.Item(122, 11).Value = 1010000.52
$ws.Cells.Item(122, 13).Value = -1007550.52
$ws.Cells.Item(126, 8).Value = 48571.285
$ws.Cells.Item(126, 9).Value = 54999.832
$ws.Cells.Item(126, 11).Value = 164999.496
$ws.Cells.Item(126, 13).Value = -162529.496
$ws.Cells.Item(132, 8).Value = 3023.6667
$ws.Cells.Item(132, 9).Value = 3121.9092
$ws.Cells.Item(132, 11).Value = 9365.7276
$ws.Cells.Item(132, 13).Value = -6835.7276

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 26318008
$ws.Cells.Item(107, 9).Value = 1747.3334
$ws.Cells.Item(107, 11).Value = 5242.0002
$ws.Cells.Item(107, 13).Value = -3322.0002
$ws.Cells.Item(132, 8).Value = 1295.5526
$ws.Cells.Item(132, 9).Value = 1139.3667
$ws.Cells.Item(132, 10).Value = 1881.25
$ws.Cells.Item(132, 11).Value = 3418.1001
$ws.Cells.Item(132, 12).Value = 5643.75
$ws.Cells.Item(132, 13).Value = -888.1001000000001
$ws.Cells.Item(132, 14).Value = -10703.75
$ws.Cells.Item(136, 8).Value = 4316.8096
$ws.Cells.Item(136, 9).Value = 3994.2727
$ws.Cells.Item(136, 10).Value = 5499.4443
$ws.Cells.Item(136, 11).Value = 11982.8181
$ws.Cells.Item(136, 12).Value = 16498.3329
$ws.Cells.Item(136, 13).Value = -9432.8181
$ws.Cells.Item(136, 14).Value = -21598.3329
